# fix 9mm damage 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("J5").Value = 30

# Row 12
$ws.Range("I12").Value = 0.06

# Row 14
$ws.Range("I14").Value = 0.14

# Row 15
$ws.Range("I15").Value = 0.14

# Row 17
$ws.Range("I17").Value = 0.14

# Row 18
$ws.Range("E18").Value = -8
$ws.Range("I18").Value = 0.14

# Row 19
$ws.Range("E19").Value = -8
$ws.Range("I19").Value = 0.14

# Row 20
$ws.Range("C20").Value = -8
$ws.Range("E20").Value = -6
$ws.Range("I20").Value = 0.14
$ws.Range("J20").Value = 200

# Row 25
$ws.Range("I25").Value = 0.05

# Update the active cell selection to match the saved view state
[void]$ws.Range("E19").Select()
